$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 1 (the "Rat/Day" label row) - shifts all data rows up by 1
$ws.Rows.Item(1).Delete()

# Delete column A (now-empty leading column) - shifts all data columns left by 1
$ws.Columns.Item(1).Delete()

# Restore the prior selection/active-cell state
$ws.Range("F22:F23").Select()
